$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Update period headers (row 8, columns D:H) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Update publish dates (row 9, columns D:H) ---
$ws.Range("D9").Value = "1399-02-01 (9)"
$ws.Range("E9").Value = "1400-02-01 (8)"
$ws.Range("F9").Value = "1401-02-07 (9)"
$ws.Range("G9").Value = "1402-02-06 (9)"
$ws.Range("H9").Value = "1402-02-06 (2)"

# --- Update financial data rows: shift one period to the left and append the new period ---
# Row 11: فروش (Sales)
$ws.Range("D11").Value = 6443
$ws.Range("E11").Value = 7077
$ws.Range("F11").Value = 8134
$ws.Range("G11").Value = 13101
$ws.Range("H11").Value = 15632

# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
$ws.Range("D12").Value = -3572
$ws.Range("E12").Value = -4021
$ws.Range("F12").Value = -3421
$ws.Range("G12").Value = -6141
$ws.Range("H12").Value = -5866

# Row 13: سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 2871
$ws.Range("E13").Value = 3056
$ws.Range("F13").Value = 4713
$ws.Range("G13").Value = 6960
$ws.Range("H13").Value = 9766

# Row 14: هزینه های عمومی, اداری و تشکیلاتی
$ws.Range("D14").Value = -522
$ws.Range("E14").Value = -539
$ws.Range("F14").Value = -427
$ws.Range("G14").Value = -568
$ws.Range("H14").Value = -657

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = -75
$ws.Range("E16").Value = -10
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -44

# Row 17: سود (زیان) عملیاتی
$ws.Range("D17").Value = 2274
$ws.Range("E17").Value = 2506
$ws.Range("F17").Value = 4289
$ws.Range("G17").Value = 6400
$ws.Range("H17").Value = 9064

# Row 18: هزینه های مالی
$ws.Range("D18").Value = -545
$ws.Range("E18").Value = -99
$ws.Range("F18").Value = -158
$ws.Range("G18").Value = -192
$ws.Range("H18").Value = -178

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 396
$ws.Range("E19").Value = 107
$ws.Range("F19").Value = 184
$ws.Range("G19").Value = 316
$ws.Range("H19").Value = 238

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 2124
$ws.Range("E20").Value = 2515
$ws.Range("F20").Value = 4315
$ws.Range("G20").Value = 6524
$ws.Range("H20").Value = 9125

# Row 21: مالیات
$ws.Range("D21").Value = -631
$ws.Range("E21").Value = -599
$ws.Range("F21").Value = -959
$ws.Range("G21").Value = -1187
$ws.Range("H21").Value = -1298

# Row 22: سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 1493
$ws.Range("E22").Value = 1916
$ws.Range("F22").Value = 3355
$ws.Range("G22").Value = 5337
$ws.Range("H22").Value = 7827

# Row 24: سود (زیان) خالص
$ws.Range("D24").Value = 1493
$ws.Range("E24").Value = 1916
$ws.Range("F24").Value = 3355
$ws.Range("G24").Value = 5337
$ws.Range("H24").Value = 7827

# Row 26: سرمایه
$ws.Range("D26").Value = 1977
$ws.Range("E26").Value = 1559
$ws.Range("F26").Value = 885
$ws.Range("G26").Value = 1137
$ws.Range("H26").Value = 850

$wb.Save()
